$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mtskheta")

# ---------------------------------------------------------------
# Insert a new row at 5 (between the old row4 "disability persons"
# row and the old row5 "Source" row) so we end up with 6 rows total:
#   1 title / 2 subtitle / 3 years / 4 family-with-disabilities /
#   5 disabilities / 6 source
# ---------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------
# Row 1 - title (merged A1:I1)
# ---------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Mtskheta Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 51

# ---------------------------------------------------------------
# Row 2 - subtitle (unchanged text/style, just row height reset)
# ---------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 14.5

# ---------------------------------------------------------------
# Row 3 - year header row; A3 switches font to Sylfaen 11
# ---------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11

# ---------------------------------------------------------------
# Row 4 - "family with disabilities Persons" (was "Number of disability persons")
# ---------------------------------------------------------------
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A4").Borders.Item(9).LineStyle = 0

$ws.Range("B4").Value = 916
$ws.Range("C4").Value = 900
$ws.Range("D4").Value = 848
$ws.Range("E4").Value = 867
$ws.Range("F4").Value = 854
$ws.Range("G4").Value = 861
$ws.Range("H4").Value = 883
$ws.Range("I4").Value = 888

$ws.Range("B4:C4").HorizontalAlignment = -4131
$ws.Range("D4:I4").Font.ColorIndex = -4105
$ws.Range("I4").Borders.Item(9).LineStyle = 0
$ws.Range("I4").Borders.Item(8).LineStyle = 0

$ws.Rows.Item(4).RowHeight = 24.75

# ---------------------------------------------------------------
# Row 5 - new "disabilities Persons" row
# ---------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.ColorIndex = -4105
$ws.Range("A5").Font.ThemeColor = 1
$ws.Range("A5").Interior.Pattern = 1
$ws.Range("A5").Interior.ThemeColor = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true

$ws.Range("B5").Value = 1034
$ws.Range("C5").Value = 1017
$ws.Range("D5").Value = 951
$ws.Range("E5").Value = 976
$ws.Range("F5").Value = 963
$ws.Range("G5").Value = 966
$ws.Range("H5").Value = 991
$ws.Range("I5").Value = 995

$ws.Range("B5:I5").NumberFormat = "#\ ##0"
$ws.Range("B5:I5").Font.Name = "Arial"
$ws.Range("B5:I5").Font.Size = 10
$ws.Range("B5:I5").Interior.Pattern = 1
$ws.Range("B5:I5").Interior.ThemeColor = 0
$ws.Range("B5:C5").Font.ColorIndex = -4105
$ws.Range("D5:I5").HorizontalAlignment = -4152
$ws.Range("I5").Borders.Item(9).LineStyle = 1

$ws.Rows.Item(5).RowHeight = 21

# ---------------------------------------------------------------
# Row 6 - Source row (previously row 5)
# ---------------------------------------------------------------
$ws.Range("A6").Value = "Source: Ministry of Internally Displaced Persons from the Occupied Territories, Labour, Health and Social Affairs of Georgia."
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Rows.Item(6).RowHeight = 27.75

# ---------------------------------------------------------------
# Column A width + default row height
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.81640625
$ws.Rows.Item(1).EntireRow.AutoFit() | Out-Null
